# Insert a new price record row for "Ají" (Cristal variety) at row 180 of
# the weekly price sheet for "Macroferia Regional de Talca". This shifts
# all the existing rows from 180..221 down by one (to 181..222) and grows
# the used range from A1:R221 to A1:R222.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 180 (pushes 180..221 -> 181..222)
$ws.Range("A180").EntireRow.Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A180").Value = 5
$ws.Range("B180").Value = "Macroferia Regional de Talca"
$ws.Range("C180").Value = "Maule"
$ws.Range("D180").Value = 44641
$ws.Range("E180").Value = 7
$ws.Range("F180").Value = 100112021
$ws.Range("G180").Value = "Ají"
$ws.Range("H180").Value = "Cristal"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 150
$ws.Range("K180").Value = 15000
$ws.Range("L180").Value = 15000
$ws.Range("M180").Value = 15000
$ws.Range("N180").Value = "`$/saco 25 kilos"
$ws.Range("O180").Value = "Región del Maule"
$ws.Range("P180").Value = 600
$ws.Range("Q180").Value = 25
$ws.Range("R180").Value = "Hortaliza"
